{"js": "// Ordered list of (oldText, newText) pairs, matching the document order of\n// the date heading followed by the 25 table-cell answers.\nconst replacements = [\n  [\"2023-09-10 Sunday\", \"2023-09-11 Monday\"],\n  [\"70\u00f72=35, 0\", \"33\u00f73=11, 0\"],\n  [\"10\u00f74=2, 2\", \"32\u00f74=8, 0\"],\n  [\"96\u00f75=19, 1\", \"17\u00f75=3, 2\"],\n  [\"32\u00f77=4, 4\", \"31\u00f79=3, 4\"],\n  [\"32\u00f78=4, 0\", \"90\u00f72=45, 0\"],\n  [\"93\u00f72=46, 1\", \"85\u00f74=21, 1\"],\n  [\"73\u00f78=9, 1\", \"31\u00f75=6, 1\"],\n  [\"30\u00f73=10, 0\", \"61\u00f72=30, 1\"],\n  [\"47\u00f75=9, 2\", \"26\u00f79=2, 8\"],\n  [\"15\u00f78=1, 7\", \"98\u00f79=10, 8\"],\n  [\"39\u00f79=4, 3\", \"99\u00f74=24, 3\"],\n  [\"28\u00f79=3, 1\", \"27\u00f75=5, 2\"],\n  [\"88\u00f79=9, 7\", \"97\u00f78=12, 1\"],\n  [\"87\u00f73=29, 0\", \"16\u00f77=2, 2\"],\n  [\"15\u00f73=5, 0\", \"46\u00f72=23, 0\"],\n  [\"18\u00f78=2, 2\", \"22\u00f73=7, 1\"],\n  [\"16\u00f74=4, 0\", \"84\u00f72=42, 0\"],\n  [\"34\u00f77=4, 6\", \"62\u00f73=20, 2\"],\n  [\"41\u00f79=4, 5\", \"35\u00f73=11, 2\"],\n  [\"74\u00f77=10, 4\", \"43\u00f73=14, 1\"],\n  [\"48\u00f78=6, 0\", \"70\u00f72=35, 0\"],\n  [\"31\u00f74=7, 3\", \"49\u00f73=16, 1\"],\n  [\"20\u00f75=4, 0\", \"56\u00f75=11, 1\"],\n  [\"19\u00f78=2, 3\", \"85\u00f77=12, 1\"],\n  [\"91\u00f77=13, 0\", \"30\u00f74=7, 2\"],\n];\n\nconst body = context.document.body;\n\n// Apply replacements one at a time (searching fresh each time) so that a\n// newly-inserted value is never accidentally matched by a later search.\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2023-09-10 Sunday\", \"2023-09-11 Monday\"),\n    @(\"70\u00f72=35, 0\", \"33\u00f73=11, 0\"),\n    @(\"10\u00f74=2, 2\", \"32\u00f74=8, 0\"),\n    @(\"96\u00f75=19, 1\", \"17\u00f75=3, 2\"),\n    @(\"32\u00f77=4, 4\", \"31\u00f79=3, 4\"),\n    @(\"32\u00f78=4, 0\", \"90\u00f72=45, 0\"),\n    @(\"93\u00f72=46, 1\", \"85\u00f74=21, 1\"),\n    @(\"73\u00f78=9, 1\", \"31\u00f75=6, 1\"),\n    @(\"30\u00f73=10, 0\", \"61\u00f72=30, 1\"),\n    @(\"47\u00f75=9, 2\", \"26\u00f79=2, 8\"),\n    @(\"15\u00f78=1, 7\", \"98\u00f79=10, 8\"),\n    @(\"39\u00f79=4, 3\", \"99\u00f74=24, 3\"),\n    @(\"28\u00f79=3, 1\", \"27\u00f75=5, 2\"),\n    @(\"88\u00f79=9, 7\", \"97\u00f78=12, 1\"),\n    @(\"87\u00f73=29, 0\", \"16\u00f77=2, 2\"),\n    @(\"15\u00f73=5, 0\", \"46\u00f72=23, 0\"),\n    @(\"18\u00f78=2, 2\", \"22\u00f73=7, 1\"),\n    @(\"16\u00f74=4, 0\", \"84\u00f72=42, 0\"),\n    @(\"34\u00f77=4, 6\", \"62\u00f73=20, 2\"),\n    @(\"41\u00f79=4, 5\", \"35\u00f73=11, 2\"),\n    @(\"74\u00f77=10, 4\", \"43\u00f73=14, 1\"),\n    @(\"48\u00f78=6, 0\", \"70\u00f72=35, 0\"),\n    @(\"31\u00f74=7, 3\", \"49\u00f73=16, 1\"),\n    @(\"20\u00f75=4, 0\", \"56\u00f75=11, 1\"),\n    @(\"19\u00f78=2, 3\", \"85\u00f77=12, 1\"),\n    @(\"91\u00f77=13, 0\", \"30\u00f74=7, 2\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n}\n"}
